# Update "想去人数" (F) and "最低票价" (G) figures across the refreshed
# data snapshot. The same underlying event data is duplicated across the
# "展览", "本地生活" and "全部类型" sheets, so each updated row is applied
# to every sheet that contains it.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) sheet ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 418
$ws1.Range("G3").Value = 60
$ws1.Range("F4").Value = 9701
$ws1.Range("F8").Value = 6617
$ws1.Range("F10").Value = 10525
$ws1.Range("F11").Value = 11694
$ws1.Range("F12").Value = 1263
$ws1.Range("F13").Value = 1215
$ws1.Range("F14").Value = 5071
$ws1.Range("F15").Value = 842
$ws1.Range("F16").Value = 506
$ws1.Range("G16").Value = 128
$ws1.Range("F20").Value = 1372
$ws1.Range("F21").Value = 286
$ws1.Range("F22").Value = 1921
$ws1.Range("F24").Value = 1334
$ws1.Range("F25").Value = 865
$ws1.Range("F27").Value = 2083
$ws1.Range("F29").Value = 674
$ws1.Range("F30").Value = 2768
$ws1.Range("F32").Value = 1853
$ws1.Range("F36").Value = 944
$ws1.Range("F37").Value = 40
$ws1.Range("F38").Value = 61
$ws1.Range("F39").Value = 3454
$ws1.Range("F40").Value = 244
$ws1.Range("F42").Value = 537
$ws1.Range("F43").Value = 606
$ws1.Range("F44").Value = 33
$ws1.Range("F45").Value = 908
$ws1.Range("F46").Value = 254
$ws1.Range("F47").Value = 14
$ws1.Range("F48").Value = 4247
$ws1.Range("F49").Value = 100

# ---- 本地生活 (Local life) sheet ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6131

# ---- 全部类型 (All types) sheet ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 418
$ws4.Range("G3").Value = 60
$ws4.Range("F4").Value = 9701
$ws4.Range("F9").Value = 10525
$ws4.Range("F10").Value = 11694
$ws4.Range("F12").Value = 1215
$ws4.Range("F13").Value = 5071
$ws4.Range("F14").Value = 842
$ws4.Range("F15").Value = 506
$ws4.Range("G15").Value = 128
$ws4.Range("F20").Value = 1372
$ws4.Range("F21").Value = 286
$ws4.Range("F22").Value = 1921
$ws4.Range("F24").Value = 1334
$ws4.Range("F25").Value = 865
$ws4.Range("F26").Value = 2083
$ws4.Range("F28").Value = 674
$ws4.Range("F29").Value = 2768
$ws4.Range("F31").Value = 1853
$ws4.Range("F39").Value = 944
$ws4.Range("F40").Value = 40
$ws4.Range("F42").Value = 244
$ws4.Range("F44").Value = 537
$ws4.Range("F45").Value = 606
$ws4.Range("F46").Value = 908
$ws4.Range("F47").Value = 254
$ws4.Range("F48").Value = 4247
$ws4.Range("F49").Value = 100
